$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# --- Updates to existing rows (new "isValidKeyReq" column) ---
$ws.Range("C2").Value = "isValidKeyReq"
$ws.Range("C3").Value = "N"

$ws.Range("F14").Value = "isValidKeyReq"
$ws.Range("F15").Value = "y"

# --- New test block: validateDeleteCustomerAPIOfValidCustomer ---
$ws.Range("A17").Value = "validateDeleteCustomerAPIOfValidCustomer"

$ws.Range("A18").Value = "endPoint"
$ws.Range("B18").Value = "methodType"
$ws.Range("C18").Value = "customerId"

$ws.Range("A19").Value = "/customers"
$ws.Range("B19").Value = "delete"
$ws.Range("C19").Value = "cus_D8uIyZJZJ2yjas"

$ws.Range("A20").Value = "endOfTestData"

# --- New test block: validateRetriveCustomersWithInvalidKey ---
$ws.Range("A21").Value = "validateRetriveCustomersWithInvalidKey"

$ws.Range("A22").Value = "endPoint"
$ws.Range("B22").Value = "expectedStatusCode"
$ws.Range("C22").Value = "isValidKeyReq"
$ws.Range("D22").Value = "methodType"
$ws.Range("E22").Value = "expectedDataSize"
$ws.Range("F22").Value = "id"

$ws.Range("A23").Value = "/customers"
$ws.Range("B23").Value = "'200"
$ws.Range("C23").Value = "y"
$ws.Range("D23").Value = "get"
$ws.Range("E23").Value = "'10"
$ws.Range("F23").Value = "cus_D9ihuaPWsiQw3K"

$ws.Range("A24").Value = "endOfTestData"

# --- New test block: verifyCustomerUsingPOJO ---
$ws.Range("A25").Value = "verifyCustomerUsingPOJO"

$ws.Range("A26").Value = "endPoint"
$ws.Range("B26").Value = "expectedStatusCode"
$ws.Range("C26").Value = "isValidKeyReq"
$ws.Range("D26").Value = "methodType"
$ws.Range("E26").Value = "expectedDataSize"
$ws.Range("F26").Value = "id"

$ws.Range("A27").Value = "/customers"
$ws.Range("B27").Value = "'200"
$ws.Range("C27").Value = "y"
$ws.Range("D27").Value = "get"
$ws.Range("E27").Value = "'10"
$ws.Range("F27").Value = "cus_D9ihuaPWsiQw3K"

$ws.Range("A28").Value = "endOfTestData"

# --- New font style (green Consolas 8pt) for the id cells ---
$ws.Range("F23").Font.Name = "Consolas"
$ws.Range("F23").Font.Size = 8
$ws.Range("F23").Font.Color = 32768

$ws.Range("F27").Font.Name = "Consolas"
$ws.Range("F27").Font.Size = 8
$ws.Range("F27").Font.Color = 32768

# --- Column widths for new columns E/F ---
$ws.Columns.Item(5).ColumnWidth = 21.88671875
$ws.Columns.Item(6).ColumnWidth = 14.88671875

# --- Selection moves to D11 ---
$ws.Activate()
$ws.Range("D11").Select()
